$wb = $excel.ActiveWorkbook

# ShipmentInformation: pickup/dropoff company name for the new claim's shipment
$wsShip = $wb.Worksheets.Item("ShipmentInformation")
$wsShip.Range("C2").NumberFormat = "@"
$wsShip.Range("C2").Value = "PickUp6899"
$wsShip.Range("K2").NumberFormat = "@"
$wsShip.Range("K2").Value = "DropOff4424"

# Input: replace the claim's order info (row 2 and row 3) with the new claim's data
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("B2").NumberFormat = "@"
$wsInput.Range("B2").Value = "07-06-2022"
$wsInput.Range("T2").NumberFormat = "@"
$wsInput.Range("T2").Value = "59082131"
$wsInput.Range("U2").NumberFormat = "@"
$wsInput.Range("U2").Value = "$697.36"
$wsInput.Range("W2").NumberFormat = "@"
$wsInput.Range("W2").Value = "999U813862"
$wsInput.Range("X2").NumberFormat = "@"
$wsInput.Range("X2").Value = "FCPBID1034723"

$wsInput.Range("B3").NumberFormat = "@"
$wsInput.Range("B3").Value = "07-06-2022"
$wsInput.Range("T3").NumberFormat = "@"
$wsInput.Range("T3").Value = "59082134"
$wsInput.Range("U3").NumberFormat = "@"
$wsInput.Range("U3").Value = "$666.00"
$wsInput.Range("W3").NumberFormat = "@"
$wsInput.Range("W3").Value = "999U814054"
$wsInput.Range("X3").NumberFormat = "@"
$wsInput.Range("X3").Value = "FCPBID1034724"

# ClaimDetail: matching claim id/date/status for the new claim
$wsClaim = $wb.Worksheets.Item("ClaimDetail")
$wsClaim.Range("A2").NumberFormat = "@"
$wsClaim.Range("A2").Value = "59082131"
$wsClaim.Range("B2").NumberFormat = "@"
$wsClaim.Range("B2").Value = "07-06-2022"
$wsClaim.Range("C2").NumberFormat = "@"
$wsClaim.Range("C2").Value = "Filed"

$wsClaim.Range("A3").NumberFormat = "@"
$wsClaim.Range("A3").Value = "59082134"
$wsClaim.Range("B3").NumberFormat = "@"
$wsClaim.Range("B3").Value = "07-06-2022"
